$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1). D1/E1/F1/G1 (eff_cap, applied voltage,
# J_ph, J) are unchanged. C1 is updated first so that the newly
# introduced shared strings end up appended in the same order as the
# authored workbook (angular frequency, Z_real, Z_imag).
$ws.Range("C1").Value = "angular frequency"
$ws.Range("A1").Value = "Z_real"
$ws.Range("B1").Value = "Z_imag"

# Update the active selection shown when the sheet is opened.
$ws.Range("B1").Select() | Out-Null
